$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update workbook/sheet title and header text for the new date
$wb.Worksheets.Item(1).Name = "Through 2022-03-04"
$ws.Range("B1").Value = "March 2022 (through March 04)"

# Neighborhood ranking shifted: Austin moved up to row 3, displacing North Lawndale and Garfield Park
$ws.Range("A3").Value = "Austin"
$ws.Range("A4").Value = "North Lawndale"
$ws.Range("A5").Value = "Garfield Park"

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 15
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 3
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 8
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 11
$ws.Range("T3").Value = 2
$ws.Range("U3").ClearContents()
$ws.Range("V3").Value = 3
$ws.Range("X3").Value = 2
$ws.Range("Y3").Value = 4

# Row 4
$ws.Range("B4").Value = 2
$ws.Range("D4").Value = 6
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 4
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 4
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 4
$ws.Range("W4").ClearContents()
$ws.Range("X4").ClearContents()
$ws.Range("Y4").Value = 3

# Row 5
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 4
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 12
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 10
$ws.Range("K5").ClearContents()
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 6
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 5
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 6
$ws.Range("W5").Value = 1
$ws.Range("X5").Value = 1
$ws.Range("Y5").Value = 1

# Row 7
$ws.Range("E7").Value = 1

# Row 16
$ws.Range("H16").Value = 2

# Row 17
$ws.Range("B17").Value = 1
$ws.Range("E17").Value = 1

# Row 18
$ws.Range("B18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("N18").Value = 1

# Row 26
$ws.Range("H26").Value = 1
$ws.Range("N26").Value = 1
$ws.Range("T26").Value = 1

# Row 32
$ws.Range("B32").Value = 1

# Row 61
$ws.Range("E61").Value = 1
$ws.Range("W61").Value = 1

# Row 75
$ws.Range("N75").Value = 1
$ws.Range("Q75").Value = 1
